$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two duplicate "UAW .../ Covestro 1-29 FEB" rows (originally rows 34 and 36)
$ws.Rows("34").Delete()
$ws.Rows("35").Delete()

# Delete the trailing "J99-9999 / TEST EXCEL" sample rows (now at 117-119)
$ws.Rows("117:119").Delete()

# Copy formatting of row 116 (A:C) down into the two new rows 117-118
$ws.Range("A116:C116").Copy($ws.Range("A117:C117"))
$ws.Range("A116:C116").Copy($ws.Range("A118:C118"))

$ws.Range("B117").Value = "Duplicate 1"
$ws.Range("B118").Value = "Duplicate 2"

$ws.Range("B107").Select()
